# Generate Report for Handback
# Removes the "1f1040ba-0717-43ae-a2e3-60273562994b" handback row (row 3) from
# every sheet, and refreshes the "Correspond Handoff/Handback DateTime"
# timestamps for the remaining "0b21aba7-7f3a-424d-92eb-a1266b37238e" row on
# the per-language sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)

# Drop the now-stale row for the 1f1040ba handback.
$wsOverview.Rows.Item(3).Delete()

# Rebuild the hyperlink collection so only the surviving row's link remains
# (row-level hyperlink deletion isn't tracked individually by this host, so
# clear everything on the sheet and re-create the one we still need).
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/52558f71b0831ee63f944a631b310650f632745f/e2e/0b21aba7-7f3a-424d-92eb-a1266b37238e.md", "", "", "0b21aba7-7f3a-424d-92eb-a1266b37238e.md")

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item(2)

# Refresh the handoff/handback timestamps for the row that survives.
$wsZhCn.Range("E2").Value = "2016-03-11 10:43:06"
$wsZhCn.Range("H2").Value = "2016-03-11 10:43:23"

# Drop the now-stale row for the 1f1040ba handback.
$wsZhCn.Rows.Item(3).Delete()

# Rebuild hyperlinks, keeping only the ones that belong to row 2.
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/52558f71b0831ee63f944a631b310650f632745f/e2e/0b21aba7-7f3a-424d-92eb-a1266b37238e.md", "", "", "0b21aba7-7f3a-424d-92eb-a1266b37238e.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/52558f71b0831ee63f944a631b310650f632745f/e2e/0b21aba7-7f3a-424d-92eb-a1266b37238e.md", "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/090c3a7dda7bc2ddabdba2c87e0fa7174a21a945/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0b21aba7-7f3a-424d-92eb-a1266b37238e.e66bd2204af96e7035b7422a7d9faa370150cab0.zh-cn.xlf", "", "", "0b21aba7-7f3a-424d-92eb-a1266b37238e.e66bd2204af96e7035b7422a7d9faa370150cab0.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/a90b72c19dd4ce253ec57617eb26d3d3b4fbccbc/e2e/0b21aba7-7f3a-424d-92eb-a1266b37238e.md", "", "", "0b21aba7-7f3a-424d-92eb-a1266b37238e.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ea6d375836f1898e4a1d69259154e9f9b5d08a41/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0b21aba7-7f3a-424d-92eb-a1266b37238e.e66bd2204af96e7035b7422a7d9faa370150cab0.zh-cn.xlf", "", "", "0b21aba7-7f3a-424d-92eb-a1266b37238e.e66bd2204af96e7035b7422a7d9faa370150cab0.zh-cn.xlf")

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item(3)

# Refresh the handoff/handback timestamps for the row that survives.
$wsDeDe.Range("E2").Value = "2016-03-11 10:43:09"
$wsDeDe.Range("H2").Value = "2016-03-11 10:43:28"

# Drop the now-stale row for the 1f1040ba handback.
$wsDeDe.Rows.Item(3).Delete()

# Rebuild hyperlinks, keeping only the ones that belong to row 2.
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/52558f71b0831ee63f944a631b310650f632745f/e2e/0b21aba7-7f3a-424d-92eb-a1266b37238e.md", "", "", "0b21aba7-7f3a-424d-92eb-a1266b37238e.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/52558f71b0831ee63f944a631b310650f632745f/e2e/0b21aba7-7f3a-424d-92eb-a1266b37238e.md", "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a2a22c60fcb6d823d6509fd0f3b66144875a8321/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0b21aba7-7f3a-424d-92eb-a1266b37238e.e66bd2204af96e7035b7422a7d9faa370150cab0.de-de.xlf", "", "", "0b21aba7-7f3a-424d-92eb-a1266b37238e.e66bd2204af96e7035b7422a7d9faa370150cab0.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/76e75520ce624216a0de609ed592893cdc3ce8ed/e2e/0b21aba7-7f3a-424d-92eb-a1266b37238e.md", "", "", "0b21aba7-7f3a-424d-92eb-a1266b37238e.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7308b6b1d69dc9aeb6e651227a39ba9d70fec970/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0b21aba7-7f3a-424d-92eb-a1266b37238e.e66bd2204af96e7035b7422a7d9faa370150cab0.de-de.xlf", "", "", "0b21aba7-7f3a-424d-92eb-a1266b37238e.e66bd2204af96e7035b7422a7d9faa370150cab0.de-de.xlf")
